# Live trading results update - Trade #22 closed.
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet: refresh headline metrics after the new closed trade
# ---------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.66   # Current Capital
$summary.Range("B4").Value = -0.34    # Total P&L $
$summary.Range("B5").Value = -0.31    # Total P&L %
$summary.Range("B6").Value = 22       # Total Trades
$summary.Range("B8").Value = 10       # Losing Trades
$summary.Range("B9").Value = 27.27    # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet: MarketMaking row (row 4) mirrors the same refresh
# ---------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.66     # Capital
$status.Range("D4").Value = 22        # Trades
$status.Range("E4").Value = -0.34     # P&L $
$status.Range("F4").Value = -0.34     # P&L %
$status.Range("G4").Value = 27.27     # Win Rate %

# ---------------------------------------------------------------
# Helper to append the new trade row (#22) to a trade log sheet
# ---------------------------------------------------------------
function Add-Trade22Row($sheet) {
    $sheet.Cells.Item(23, 1).Value = 22
    # Date column stores a plain text date string (not a real Excel date) in
    # this workbook, so force text formatting before assignment to avoid the
    # COM layer auto-converting the string into a date serial number, then
    # reset the style back to Normal so no stray formatting is introduced.
    $dateCell = $sheet.Cells.Item(23, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"
    $sheet.Cells.Item(23, 3).Value = "08:02:43"
    $sheet.Cells.Item(23, 4).Value = "MarketMaking"
    $sheet.Cells.Item(23, 5).Value = "DOWN"
    $sheet.Cells.Item(23, 6).Value = 0.57
    $sheet.Cells.Item(23, 7).Value = 0.41
    $sheet.Cells.Item(23, 8).Value = "CLOSED"
    $sheet.Cells.Item(23, 9).Value = -28.0702
    $sheet.Cells.Item(23, 10).Value = -0.16
    $sheet.Cells.Item(23, 11).Value = 99.66
    $sheet.Cells.Item(23, 12).Value = 0
    $sheet.Cells.Item(23, 13).Value = 0
    $sheet.Cells.Item(23, 14).Value = 0.6
    $sheet.Cells.Item(23, 15).Value = "Normal spread capture: 19600 bps"
    $sheet.Cells.Item(23, 16).Value = "early_exit"
    $sheet.Cells.Item(23, 17).Value = 0.13
}

# ---------------------------------------------------------------
# All Trades sheet: append trade #22
# ---------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade22Row $allTrades

# ---------------------------------------------------------------
# MarketMaking sheet: append trade #22 (same log, strategy-specific view)
# ---------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade22Row $marketMaking
